$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.657.76'
$ws.Range("E2").Value = '  +2.64%  '
$ws.Range("D3").Value = '1.788.07'
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '223.12'
$ws.Range("E5").Value = '  -0.97%  '
$ws.Range("E6").Value = '  -0.76%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '32.63'
$ws.Range("E8").Value = '  +7.03%  '
$ws.Range("E9").Value = '  +0.60%  '
$ws.Range("D10").Value = '0.0680'
$ws.Range("E10").Value = '  +2.04%  '
$ws.Range("E11").Value = '  +1.50%  '
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("D13").Value = '11.18'
$ws.Range("E13").Value = '  +11.30%  '
$ws.Range("D14").Value = '1.788.92'
$ws.Range("E14").Value = '  +0.78%  '
$ws.Range("D15").Value = '34.628.98'
$ws.Range("E15").Value = '  +2.61%  '
$ws.Range("E16").Value = '  +0.96%  '
$ws.Range("D17").Value = '4.32'
$ws.Range("E17").Value = '  +3.24%  '
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").Value = '253.32'
$ws.Range("E19").Value = '  +0.56%  '
$ws.Range("D20").Value = '0.0₃0774'
$ws.Range("E20").Value = '  +4.89%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("E22").Value = '  +1.57%  '
$ws.Range("E23").Value = '  +0.71%  '
$ws.Range("D24").Value = '2.13'
$ws.Range("E24").Value = '  -0.37%  '
$ws.Range("D25").Value = '158.47'
$ws.Range("E25").Value = '  -0.58%  '
$ws.Range("D26").Value = '16.33'
$ws.Range("E26").Value = '  -0.93%  '
$ws.Range("E27").Value = '  +1.97%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").Value = '3.76'
$ws.Range("E30").Value = '  -1.48%  '
$ws.Range("D31").Value = '0.0514'
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("E32").Value = '  -0.38%  '
$ws.Range("E33").Value = '  +0.42%  '
$ws.Range("E34").Value = '  -1.39%  '
$ws.Range("D35").Value = '1.439.74'
$ws.Range("E35").Value = '  -2.81%  '
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("D37").Value = '0.0190'
$ws.Range("E37").Value = '  +2.65%  '
$ws.Range("E38").Value = '  -0.25%  '
$ws.Range("D39").Value = '83.04'
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("E40").Value = '  +4.44%  '
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("D42").Value = '0.905'
$ws.Range("E42").Value = '  +2.33%  '
$ws.Range("E43").Value = '  -0.67%  '
$ws.Range("D44").Value = '0.0507'
$ws.Range("E44").Value = '  -0.60%  '
$ws.Range("D45").Value = '5.96'
$ws.Range("E45").Value = '  +4.17%  '
$ws.Range("E46").Value = '  -2.60%  '
$ws.Range("D47").Value = '1.944.18'
$ws.Range("E47").Value = '  +0.63%  '
$ws.Range("D48").Value = '105.11'
$ws.Range("E48").Value = '  +7.76%  '
$ws.Range("D49").Value = '12.02'
$ws.Range("E49").Value = '  +1.53%  '
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = '49.54'
$ws.Range("E51").Value = '  -2.38%  '
